# "Test queries updated slightly, fixed requery bug"
#
# Fixes a typo in two of the "(987) DKK HUF ..." test-query strings:
# an extra trailing "0" in the year made the date "20012" instead of
# "2012". Also moves the sheet's active selection down to A20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the wide-spaced variant first, then the tight variant, so that
# any newly created shared-string entries land in the same append order
# as the authoritative edit.
$ws.Range("A25").Value = "   (987)     DKK     HUF           12-12-2012     "
$ws.Range("A20").Value = "(987) DKK HUF 12-12-2012"

# Move the selection/cursor to A20, scrolled so row 11 is at the top.
$ws.Range("A20").Select()
